# Populate rows 5-25 with the new Response/API question-function mapping rows.
# Each new row is built from: the Column-A question text, the Column-B API
# function name, and the row (above row 5) whose look (font/alignment) it should
# copy, per column -- since those existing rows already carry the right style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 5; A = "Return total building energy consumption"; B = "get_building_energy_consumption_overall"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 6; A = "Return total building energy consumption with a break down by end use category"; B = "get_building_energy_consumption_by_end_uses_category"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 7; A = "Return aggregated energy consumption for a specific end use category"; B = "get_building_energy_consumption_by_end_uses_category"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 8; A = "Return energy consumption for a specific piece of equipment"; B = "get_most_consumption_equipment"; AStyleFrom = "A2"; BStyleFrom = "B2" },
    @{ Row = 9; A = "Return top 5 pieces of equipment that contribute most to energy consumption of a specific end use category"; B = "get_explorer_equip_power_consumption"; AStyleFrom = "A2"; BStyleFrom = "B2" },
    @{ Row = 10; A = "Return top 5 pieces of equipment that contribute most to energy consumption for a specific equipment type"; B = "get_explorer_equip_power_consumption"; AStyleFrom = "A2"; BStyleFrom = "B2" },
    @{ Row = 11; A = "Return top 5 pieces of equipment that contribute most to energy consumption for a specific time period"; B = "get_energy_building_equipment"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 12; A = "Return total building after hours energy consumption (and % of total)"; B = "get_building_energy_consumption_by_end_uses_category"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 13; A = "Return aggregated after hours energy consumption for a specific end use category (and % of total of end use category)"; B = "get_building_energy_consumption_by_end_uses_category"; AStyleFrom = "A3"; BStyleFrom = "B3" },
    @{ Row = 14; A = "Return highest power draw of equipment in specified time period and time of occurence."; B = "get_most_consumption_equipment"; AStyleFrom = "A2"; BStyleFrom = "B2" },
    @{ Row = 15; A = "Return time of day when average energy consumption was the highest over specified period (and kWh value)"; B = "get_power_consumption_weekdays_weekend"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 16; A = "Return energy consumption by month for [total building]"; B = "get_building_energy_consumption_overall"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 17; A = "Return magintude of change for total building energy consumption from one period to the next (i.e., this month vs. last month)"; B = "get_building_energy_consumption_overall"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 18; A = "Return contribution of change to total building energy consumption broken down by [end-use category]"; B = "get_building_energy_consumption_by_end_uses_category"; AStyleFrom = "A3"; BStyleFrom = "B3" },
    @{ Row = 19; A = "Return contribution of change to total building energy consumption broken down by [specific equipment] (show top 5 in magnitude)"; B = "get_energy_building_equipment"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 20; A = "Return magnitude of change in energy consumption for [a specific end use category]"; B = "get_building_energy_consumption_by_end_uses_category"; AStyleFrom = "A3"; BStyleFrom = "B3" },
    @{ Row = 21; A = "Return magnitude of change in energy consumption for [a specific piece of equipment]"; B = "get_most_consumption_equipment"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 22; A = "Return top 5 contributors (specific pieces of equipment) that contributed to a change in energy consumption for [total building]"; B = "get_explorer_equip_power_consumption"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 23; A = "Return top 5 contributors (specific pieces of equipment) that contributed to a change in energy consumption for [a specific end-use category]"; B = "get_explorer_equip_power_consumption"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 24; A = "Return top 5 contributors (specific pieces of equipment) that contributed to a change in energy consumption for [a specific equipment type]"; B = "get_explorer_equip_power_consumption"; AStyleFrom = "A3"; BStyleFrom = "B2" },
    @{ Row = 25; A = "Return top 5 contributors (specific pieces of equipment) that contributed to a change in energy consumption for [a specific space type]"; B = "get_explorer_equip_power_consumption"; AStyleFrom = "A3"; BStyleFrom = "B2" }
)

foreach ($row in $rows) {
    $ws.Range("A$($row.Row)").Value = $row.A
    $ws.Range("B$($row.Row)").Value = $row.B

    $ws.Range($row.AStyleFrom).Copy() | Out-Null
    $ws.Range("A$($row.Row)").PasteSpecial(-4122) | Out-Null
    $ws.Range($row.BStyleFrom).Copy() | Out-Null
    $ws.Range("B$($row.Row)").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# --- Selection matches the author's last-saved cursor position ---
$ws.Range("B30").Select()
